# Weekly update: insert this week's Arandano (blue) price record at the
# top of the data table (row 5), pushing all existing rows down by one.
# The new record repeats the prior week's price figures (previous row 5,
# now row 6) but dated one week later (44544 = 2021-12-14).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 5; rows 5-13 shift down to 6-14.
$ws.Rows.Item(5).Insert()

# Populate the new row 5 with this week's record.
$ws.Cells.Item(5, 1).Value  = 4
$ws.Cells.Item(5, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(5, 3).Value  = "Los Lagos"
$ws.Cells.Item(5, 4).Value  = 44544
$ws.Cells.Item(5, 5).Value  = 10
$ws.Cells.Item(5, 6).Value  = "Fruta"
$ws.Cells.Item(5, 7).Value  = 100101
$ws.Cells.Item(5, 8).Value  = "Berries"
$ws.Cells.Item(5, 9).Value  = 100101001
$ws.Cells.Item(5, 10).Value = "Arándano (blue)"
$ws.Cells.Item(5, 11).Value = "Sin especificar"
$ws.Cells.Item(5, 12).Value = "Primera"
$ws.Cells.Item(5, 13).Value = 400
$ws.Cells.Item(5, 14).Value = 5000
$ws.Cells.Item(5, 15).Value = 5500
$ws.Cells.Item(5, 16).Value = 5250
$ws.Cells.Item(5, 17).Value = "`$/bandeja 12 canastillos 125 gramos"
$ws.Cells.Item(5, 18).Value = "Región del Maule"
$ws.Cells.Item(5, 19).Value = 3500
$ws.Cells.Item(5, 20).Value = 1.5
